$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated/new data for rows 2-8, columns A:E (epoch, train_loss, val_loss, train_mae, val_mae)
$data = @(
    @(1, 0.9334660564871232, 6.080389022827148, 6.73623922529464,  6.080389022827148),
    @(2, 0.9097319368765692, 6.098905086517334, 6.64501003820769,  6.098905086517334),
    @(3, 0.9048111958044588, 6.108983993530273, 6.625022171821506, 6.108983993530273),
    @(4, 0.9013950093163139, 6.109508991241455, 6.609887238721559, 6.109508991241455),
    @(5, 0.8988497905485038, 6.119451999664307, 6.598473477806126, 6.119451999664307),
    @(6, 0.8971853885462555, 6.11155366897583,  6.591057170571418, 6.11155366897583),
    @(7, 0.8958477337786327, 6.108453273773193, 6.585326235144985, 6.108453273773193)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Count; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $vals[$j]
    }
}
